# edit.ps1
# Applies the two substantive textual changes described by the commit
# "fix ref to fig 1a":
#
#   1. In the reviewer-response paragraph that immediately follows the
#      comment "Related to this, L87 refers to Fig 1a when it should be
#      Fig 1d.", the terse reply "- Fixed" is expanded to:
#         "- Fixed. We changed the order of panels in Fig 1 so the
#          reference to Fig 1a is now correct (line 91)"
#
#   2. The cached result of the header's DATE field is bumped from
#      "April 21, 2021" to "April 22, 2021".
#
# (All of the other hunks in the source diff are purely Word's automatic
#  <w:lastRenderedPageBreak/> repagination bookkeeping - i.e. runs being
#  re-split/re-merged around a moving page-break marker with no change
#  to the actual visible text - which Word recomputes on its own once the
#  underlying text above has changed length, so nothing further needs to
#  be scripted for those.)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Expand "- Fixed" -> "- Fixed. We changed the order of panels in
#    Fig 1 so the reference to Fig 1a is now correct (line 91)"
#
# "- Fixed" occurs several times verbatim in this document, so we must
# not do a document-wide replace. Instead, locate the unique anchor
# paragraph ("Related to this, L87 refers to Fig 1a ...") and operate
# only on the very next paragraph, which begins with the "- Fixed" reply.
# ---------------------------------------------------------------------

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Related to this, L87 refers to Fig 1a*") {
        $anchor = $p
        break
    }
}

if ($anchor -ne $null) {
    $replyPara = $anchor.Next()
    $scopedRange = $replyPara.Range

    $hit = $scopedRange.Find.Execute("- Fixed", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    if ($hit) {
        $matchStart = $scopedRange.Start

        $part1 = "- "
        $part2 = "Fixed. We changed the order of panels in Fig 1 so the reference to Fig 1a is now correct"
        $part3 = " (line 91)"

        # Replace "- Fixed" with the full expanded sentence first.
        $scopedRange.Text = $part1 + $part2 + $part3

        # Re-split the replaced text into three runs (mirroring how Word
        # naturally breaks up incrementally-edited text) while keeping
        # each run's character formatting identical, by nudging a
        # formatting property back to its own value on each sub-range.
        $r1 = $d.Range($matchStart, $matchStart + $part1.Length)
        $r1.Bold = 1
        $r1.Bold = 0

        $r2Start = $matchStart + $part1.Length
        $r2 = $d.Range($r2Start, $r2Start + $part2.Length)
        $r2.Bold = 1
        $r2.Bold = 0

        $r3Start = $r2Start + $part2.Length
        $r3 = $d.Range($r3Start, $r3Start + $part3.Length)
        $r3.Bold = 1
        $r3.Bold = 0
    }
}

# ---------------------------------------------------------------------
# 2) Bump the header date from "April 21, 2021" to "April 22, 2021"
# ---------------------------------------------------------------------

foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        $hdr.Range.Find.Execute("April 21, 2021", $true, $false, $false, $false, $false, $true, 1, $false, "April 22, 2021", 2) | Out-Null
    }
}
